$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (shifts existing rows 3..73 down to 4..74,
# carrying their formatting/styles along, same as a manual Excel row insert)
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new price-quote record
$ws.Range("A3").Value = 7
$ws.Range("B3").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C3").Value = "Ñuble"
$ws.Range("D3").Value = 44643
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 100112021
$ws.Range("G3").Value = "Ají"
$ws.Range("H3").Value = "Americana (o)"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 8500
$ws.Range("L3").Value = 9000
$ws.Range("M3").Value = 8750
$ws.Range("N3").Value = "$/caja 15 kilos"
$ws.Range("O3").Value = "Región del Maule"
$ws.Range("P3").Value = 583
$ws.Range("Q3").Value = 15
$ws.Range("R3").Value = "Hortaliza"
